$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (shashikumar) is updated in place:
#  - roll number (B2) becomes text instead of a number
#  - date (C2) and time (D2) are refreshed with a new attendance entry
# Force the target cells to text first so Excel does not reinterpret the
# date/time-looking strings as real dates/numbers, then drop the temporary
# formatting again so the cells are left without an explicit style, just
# like the rest of the data rows.
$ws.Range("B2:D2").NumberFormat = "@"
$ws.Range("B2").Value = "190001055"
$ws.Range("C2").Value = "2025-11-01"
$ws.Range("D2").Value = "10:44:25"
$ws.Range("B2:D2").ClearFormats()

# Rows 3-5 (the other attendance entries) are removed entirely, shrinking
# the sheet down to just the header row and the single updated row.
$ws.Range("A3:D5").EntireRow.Delete()
